$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 already carries the exact cell styles this new row needs
# (plain centered number, unstyled text x2, m/d/yyyy date), so clone its
# formatting down into row 56 before setting values.
$ws.Range("A16:D16").Copy()
$ws.Range("A56:D56").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(56, 1).Value = 152
$ws.Cells.Item(56, 2).Value = "Maximum Product Subarray"
$ws.Cells.Item(56, 3).Value = "Java"
$ws.Cells.Item(56, 4).Value = (Get-Date -Year 2023 -Month 5 -Day 11).Date

$ws.Range("D57").Select()
